$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. It belongs right
# before the existing row 204, so insert a blank row there first - this
# shifts the old rows 204:330 down to 205:331 (and the sheet's used range
# grows from R330 to R331 automatically).
$ws.Rows("204").Insert()

# Fill in the newly inserted row with the new record's data.
$ws.Range("A204").Value = 4
$ws.Range("B204").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C204").Value = "Los Lagos"
$ws.Range("D204").Value = 44777
$ws.Range("E204").Value = 10
$ws.Range("F204").Value = 100114014
$ws.Range("G204").Value = "Betarraga"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 500
$ws.Range("K204").Value = 1000
$ws.Range("L204").Value = 1200
$ws.Range("M204").Value = 1100
$ws.Range("N204").Value = "$/paquete 5 unidades"
$ws.Range("O204").Value = "Región del Maule"
$ws.Range("P204").Value = 220
$ws.Range("Q204").Value = 5
$ws.Range("R204").Value = "Hortaliza"
